$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Schedule": the stray capitalised "Free" is normalised to the
# lowercase "free" already used everywhere else on the sheet.
# ---------------------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")
$schedule.Activate()
$schedule.Range("F2").Value = "free"
$schedule.Range("F2").Select()

# ---------------------------------------------------------------------------
# Sheet "Week": the reader that used to mark every slot with "I" (Incomplete)
# now marks it "P" (Present) once the supply reader actually resolves a
# match -- except for three cells where a real on-caller code number was
# found, so those get written as raw numbers instead of the placeholder.
# ---------------------------------------------------------------------------
$week = $wb.Worksheets.Item("Week")
$week.Activate()

$week.Range("B2").Value = "P"
$week.Range("D2").Value = "P"
$week.Range("F2").Value = "P"

$week.Range("B3").Value = "P"
$week.Range("C3").Value = 62
$week.Range("F3").Value = "P"
$week.Range("G3").Value = "P"

$week.Range("B4").Value = "P"
$week.Range("C4").Value = "P"
$week.Range("D4").Value = 123
$week.Range("E4").Value = "P"
$week.Range("F4").Value = "P"
$week.Range("G4").Value = "P"

$week.Range("B5").Value = 151
$week.Range("C5").Value = "P"
$week.Range("D5").Value = "P"
$week.Range("E5").Value = "P"
$week.Range("F5").Value = "P"
$week.Range("G5").Value = "P"

# Selection left where the author's run ended up (Week stays the active tab).
$week.Range("G4").Select()
